$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.227.36"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.856.39"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.56"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6991"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07776"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3075"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.81"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07816"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("D12").Value = "1.856.31"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.108"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.18"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6873"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.520"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008437"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Value = "29.218.37"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.64"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "2.108.78"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.525"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9996"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.18"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.868"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.560"
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.248"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05204"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7604"
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.169"
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.844"
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.707"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "1.220.50"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.725"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8993"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.89"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9990"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.508"
$ws.Range("E44").Value = "  -11.68%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "2.007.06"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("E46").Value = "  -3.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.56"
$ws.Range("E47").Value = "  -7.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.573"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.027"
$ws.Range("E51").Value = "  +0.61%  "
